$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 73
$ws.Range("A73").Value = 22033030
$ws.Range("B73").Value = "木质修理锤"
$ws.Range("C73").Value = 1
$ws.Range("D73").Value = "立即回复王塔200生命"
$ws.Range("E73").Value = 1
$ws.Range("F73").Value = 3
$ws.Range("G73").Value = 99
$ws.Range("H73").Value = 200
$ws.Range("I73").Value = 11
$ws.Range("M73").Value = 1
$ws.Range("N73").Value = "xiulichui1"

# Row 74
$ws.Range("A74").Value = 22033031
$ws.Range("B74").Value = "钢铁修理锤"
$ws.Range("C74").Value = 1
$ws.Range("D74").Value = "立即回复王塔500生命"
$ws.Range("E74").Value = 1
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = 99
$ws.Range("H74").Value = 200
$ws.Range("I74").Value = 11
$ws.Range("M74").Value = 1
$ws.Range("N74").Value = "xiulichui2"

# Row 75
$ws.Range("A75").Value = 22033032
$ws.Range("B75").Value = "神圣修理锤"
$ws.Range("C75").Value = 1
$ws.Range("D75").Value = "立即回复王塔1000生命"
$ws.Range("E75").Value = 1
$ws.Range("F75").Value = 5
$ws.Range("G75").Value = 99
$ws.Range("H75").Value = 200
$ws.Range("I75").Value = 11
$ws.Range("M75").Value = 1
$ws.Range("N75").Value = "xiulichui3"

# Row 76
$ws.Range("A76").Value = 22034001
$ws.Range("B76").Value = "经验之书"
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = "使用后获得100点经验值"
$ws.Range("E76").Value = 1
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 99
$ws.Range("H76").Value = 200
$ws.Range("I76").Value = 13
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = "huishu"

# Row 77
$ws.Range("A77").Value = 22034002
$ws.Range("B77").Value = "能量之书"
$ws.Range("C77").Value = 1
$ws.Range("D77").Value = "使用后获得500点经验值"
$ws.Range("E77").Value = 1
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 99
$ws.Range("H77").Value = 200
$ws.Range("I77").Value = 13
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = "hongshu"

# Row 78
$ws.Range("A78").Value = 22034003
$ws.Range("B78").Value = "攻速药水"
$ws.Range("C78").Value = 1
$ws.Range("D78").Value = "提升最后召唤的生物5点攻速"
$ws.Range("E78").Value = 1
$ws.Range("F78").Value = 2
$ws.Range("G78").Value = 99
$ws.Range("H78").Value = 200
$ws.Range("I78").Value = 11
$ws.Range("M78").Value = 4
$ws.Range("N78").Value = "yaoshuistr"

# Row 79
$ws.Range("A79").Value = 22034004
$ws.Range("B79").Value = "守护药水"
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = "提升最后召唤的生物5点防御"
$ws.Range("E79").Value = 1
$ws.Range("F79").Value = 2
$ws.Range("G79").Value = 99
$ws.Range("H79").Value = 200
$ws.Range("I79").Value = 11
$ws.Range("M79").Value = 4
$ws.Range("N79").Value = "yaoshuidef"

# Row 80
$ws.Range("A80").Value = 22034005
$ws.Range("B80").Value = "法术药水"
$ws.Range("C80").Value = 1
$ws.Range("D80").Value = "提升最后召唤的生物5点魔力"
$ws.Range("E80").Value = 1
$ws.Range("F80").Value = 2
$ws.Range("G80").Value = 99
$ws.Range("H80").Value = 200
$ws.Range("I80").Value = 11
$ws.Range("M80").Value = 4
$ws.Range("N80").Value = "yaoshuimag"

# Row 81
$ws.Range("A81").Value = 22034006
$ws.Range("B81").Value = "技巧药水"
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = "提升最后召唤的生物5点命中"
$ws.Range("E81").Value = 1
$ws.Range("F81").Value = 2
$ws.Range("G81").Value = 99
$ws.Range("H81").Value = 200
$ws.Range("I81").Value = 11
$ws.Range("M81").Value = 4
$ws.Range("N81").Value = "yaoshuiskl"

# Row 82
$ws.Range("A82").Value = 22034007
$ws.Range("B82").Value = "速度药水"
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = "提升最后召唤的生物5点回避"
$ws.Range("E82").Value = 1
$ws.Range("F82").Value = 2
$ws.Range("G82").Value = 99
$ws.Range("H82").Value = 200
$ws.Range("I82").Value = 11
$ws.Range("M82").Value = 4
$ws.Range("N82").Value = "yaoshuispd"

# Row 83
$ws.Range("A83").Value = 22034008
$ws.Range("B83").Value = "幸运药水"
$ws.Range("C83").Value = 1
$ws.Range("D83").Value = "提升最后召唤的生物5点幸运"
$ws.Range("E83").Value = 1
$ws.Range("F83").Value = 2
$ws.Range("G83").Value = 99
$ws.Range("H83").Value = 200
$ws.Range("I83").Value = 11
$ws.Range("M83").Value = 4
$ws.Range("N83").Value = "yaoshuiluk"

# Row 84
$ws.Range("A84").Value = 22034009
$ws.Range("B84").Value = "暴击药水"
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = "提升最后召唤的生物5点暴击"
$ws.Range("E84").Value = 1
$ws.Range("F84").Value = 2
$ws.Range("G84").Value = 99
$ws.Range("H84").Value = 200
$ws.Range("I84").Value = 11
$ws.Range("M84").Value = 4
$ws.Range("N84").Value = "yaoshuivit"

# Row 85
$ws.Range("A85").Value = 22034010
$ws.Range("B85").Value = "饼干"
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = "使用后获得50点食物"
$ws.Range("E85").Value = 1
$ws.Range("F85").Value = 2
$ws.Range("G85").Value = 99
$ws.Range("H85").Value = 200
$ws.Range("I85").Value = 13
$ws.Range("M85").Value = 4
$ws.Range("N85").Value = "bingan"

# Row 86
$ws.Range("A86").Value = 22034011
$ws.Range("B86").Value = "红色胶囊"
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = "使用后获得50点健康"
$ws.Range("E86").Value = 1
$ws.Range("F86").Value = 2
$ws.Range("G86").Value = 99
$ws.Range("H86").Value = 200
$ws.Range("I86").Value = 13
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = "pillred"

# Row 87
$ws.Range("A87").Value = 22034012
$ws.Range("B87").Value = "蓝色胶囊"
$ws.Range("C87").Value = 1
$ws.Range("D87").Value = "使用后获得50点精神"
$ws.Range("E87").Value = 1
$ws.Range("F87").Value = 2
$ws.Range("G87").Value = 99
$ws.Range("H87").Value = 200
$ws.Range("I87").Value = 13
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = "pillblue"

# Row 88
$ws.Range("A88").Value = 22034013
$ws.Range("B88").Value = "水晶球"
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = "使用后可以看到本场景的所有事件"
$ws.Range("E88").Value = 1
$ws.Range("F88").Value = 5
$ws.Range("G88").Value = 99
$ws.Range("H88").Value = 200
$ws.Range("I88").Value = 13
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = "shuijingqiu"
